$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'66.121.30"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").Formula = "'3.502.24"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Formula = "'581.86"
$ws.Range("E5").Value = "  +5.14%  "

$ws.Range("D6").Formula = "'177.47"
$ws.Range("E6").Value = "  -6.21%  "

$ws.Range("D7").Formula = "'0.632"
$ws.Range("E7").Value = "  +3.92%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Formula = "'0.637"
$ws.Range("E9").Value = "  +0.31%  "

$ws.Range("D10").Formula = "'0.160"
$ws.Range("E10").Value = "  +3.58%  "

$ws.Range("D11").Formula = "'55.65"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").Formula = "'0.0000278"
$ws.Range("E12").Value = "  +2.47%  "

$ws.Range("D13").Formula = "'9.27"
$ws.Range("E13").Value = "  -1.53%  "

$ws.Range("D14").Formula = "'4.064.02"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").Formula = "'3.509.99"
$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Formula = "'18.29"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Formula = "'66.114.23"
$ws.Range("E18").Value = "  -1.80%  "

$ws.Range("D19").Formula = "'11.98"
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("E20").Value = "  +1.67%  "

$ws.Range("D21").Formula = "'411.42"
$ws.Range("E21").Value = "  -4.14%  "

$ws.Range("D22").Formula = "'4.27"
$ws.Range("E22").Value = "  +8.83%  "

$ws.Range("D23").Formula = "'4.41"
$ws.Range("E23").Value = "  +6.90%  "

$ws.Range("D24").Formula = "'84.87"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").Formula = "'13.34"
$ws.Range("E25").Value = "  +10.13%  "

$ws.Range("D26").Formula = "'11.04"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("D27").Formula = "'2.85"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("E28").Value = "  -1.74%  "

$ws.Range("D29").Formula = "'9.13"
$ws.Range("E29").Value = "  +1.20%  "

$ws.Range("D30").Formula = "'30.28"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").Formula = "'6.66"
$ws.Range("E31").Value = "  -0.28%  "

$ws.Range("D32").Formula = "'11.73"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").Formula = "'596.87"
$ws.Range("E33").Value = "  -7.37%  "

$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Formula = "'60.74"
$ws.Range("E35").Value = "  +1.62%  "

$ws.Range("E36").Value = "  +4.55%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +9.18%  "

$ws.Range("D39").Formula = "'" + "0.0" + [char]0x2083 + "0793"
$ws.Range("E39").Value = "  -5.10%  "

$ws.Range("D40").Formula = "'36.75"
$ws.Range("E40").Value = "  -4.97%  "

$ws.Range("D41").Formula = "'0.384"
$ws.Range("E41").Value = "  -2.05%  "

$ws.Range("D42").Formula = "'3.237.29"
$ws.Range("E42").Value = "  +6.25%  "

$ws.Range("D43").Formula = "'1.00"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").Formula = "'2.95"
$ws.Range("E44").Value = "  +2.35%  "

$ws.Range("D45").Formula = "'3.33"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Formula = "'2.54"
$ws.Range("E46").Value = "  -4.73%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Formula = "'0.0420"
$ws.Range("E47").Value = "  +0.36%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Formula = "'2.68"
$ws.Range("E48").Value = "  -5.65%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Formula = "'0.132"
$ws.Range("E49").Value = "  +1.13%  "

$ws.Range("D50").Formula = "'8.56"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("D51").Formula = "'138.23"
$ws.Range("E51").Value = "  -2.15%  "
